# Daily attendance processing - 2025-12-30 05:39:30
#
# Normalize the "Recorded By" column (G) on the "Session Analysis Results"
# sheet so that entries listing both the system and a user account show
# "System" first, e.g.:
#   "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

# Column G = "Recorded By". Use Excel's native Find & Replace against the
# whole column so only cells whose full text matches the old value are
# updated (xlWhole = 1, xlByRows = 1).
$xlWhole = 1
$xlByRows = 1

$ws.Columns.Item(7).Replace($oldValue, $newValue, $xlWhole, $xlByRows, $false, $false, $true)
